$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style for plain (unstyled) text cells, taken from an untouched D-column cell
$normalStyle = $ws.Range("D4").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.146.78"
$ws.Range("D2").Style = $normalStyle
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.053.61"
$ws.Range("D3").Style = $normalStyle
$ws.Range("E3").Value = "  -0.07%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.55"
$ws.Range("D5").Style = $normalStyle
$ws.Range("E5").Value = "  +0.45%  "

$ws.Range("E6").Value = "  +1.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.72"
$ws.Range("D7").Style = $normalStyle
$ws.Range("E7").Value = "  +5.81%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.382"
$ws.Range("D9").Style = $normalStyle
$ws.Range("E9").Value = "  +0.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0784"
$ws.Range("D10").Style = $normalStyle
$ws.Range("E10").Value = "  -1.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.09"
$ws.Range("D12").Style = $normalStyle
$ws.Range("E12").Value = "  +5.98%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.351.91"
$ws.Range("D13").Style = $normalStyle
$ws.Range("E13").Value = "  -0.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.816"
$ws.Range("D14").Style = $normalStyle
$ws.Range("E14").Value = "  -0.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.54"
$ws.Range("D15").Style = $normalStyle
$ws.Range("E15").Value = "  +5.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.054.70"
$ws.Range("D16").Style = $normalStyle
$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "37.139.33"
$ws.Range("D17").Style = $normalStyle
$ws.Range("E17").Value = "  +0.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.13"
$ws.Range("D18").Style = $normalStyle
$ws.Range("E18").Value = "  +20.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.71"
$ws.Range("D19").Style = $normalStyle
$ws.Range("E19").Value = "  +3.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0902"
$ws.Range("D20").Style = $normalStyle
$ws.Range("E20").Value = "  -0.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.37"
$ws.Range("D21").Style = $normalStyle
$ws.Range("E21").Value = "  +1.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.82"
$ws.Range("D22").Style = $normalStyle
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.19"
$ws.Range("D25").Style = $normalStyle
$ws.Range("E25").Value = "  +11.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.34"
$ws.Range("D26").Style = $normalStyle
$ws.Range("E26").Value = "  -0.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.30"
$ws.Range("D27").Style = $normalStyle
$ws.Range("E27").Value = "  +3.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.92"
$ws.Range("D28").Style = $normalStyle
$ws.Range("E28").Value = "  -0.99%  "

$ws.Range("E29").Value = "  +0.86%  "

$ws.Range("E30").Value = "  +9.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.73"
$ws.Range("D31").Style = $normalStyle
$ws.Range("E31").Value = "  +4.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0615"
$ws.Range("D32").Style = $normalStyle
$ws.Range("E32").Value = "  -1.40%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.50"
$ws.Range("D33").Style = $normalStyle
$ws.Range("E33").Value = "  +4.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0895"
$ws.Range("D34").Style = $normalStyle
$ws.Range("E34").Value = "  +4.24%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.24"
$ws.Range("D36").Style = $normalStyle
$ws.Range("E36").Value = "  -0.82%  "

$ws.Range("E37").Value = "  -2.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.109"
$ws.Range("D38").Style = $normalStyle
$ws.Range("E38").Value = "  +5.72%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.36"
$ws.Range("D39").Style = $normalStyle
$ws.Range("E39").Value = "  +0.69%  "

$ws.Range("E40").Value = "  +13.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.23"
$ws.Range("D41").Style = $normalStyle
$ws.Range("E41").Value = "  +30.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.54"
$ws.Range("D42").Style = $normalStyle
$ws.Range("E42").Value = "  -3.55%  "

$ws.Range("E43").Value = "  -0.89%  "

$ws.Range("E44").Value = "  -1.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.12"
$ws.Range("D45").Style = $normalStyle
$ws.Range("E45").Value = "  +0.18%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.45"
$ws.Range("D46").Style = $normalStyle
$ws.Range("E46").Value = "  +2.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.283.18"
$ws.Range("D47").Style = $normalStyle
$ws.Range("E47").Value = "  -1.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.87"
$ws.Range("D48").Style = $normalStyle
$ws.Range("E48").Value = "  -1.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.80"
$ws.Range("D49").Style = $normalStyle
$ws.Range("E49").Value = "  +0.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.242.03"
$ws.Range("D50").Style = $normalStyle
$ws.Range("E50").Value = "  +0.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.50"
$ws.Range("D51").Style = $normalStyle
$ws.Range("E51").Value = "  -17.13%  "
